$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 579-580; this pushes the former rows 579-598 down to 581-600
# and matches the new sheet dimension A1:R600.
$ws.Rows("579:580").Insert()

# --- New row 579 ---
$ws.Range("A579").Value = 3
$ws.Range("B579").Value = "Femacal de La Calera"
$ws.Range("C579").Value = "Coquimbo"
$ws.Range("D579").Value = 45075
$ws.Range("E579").Value = 5
$ws.Range("F579").Value = 100112031
$ws.Range("G579").Value = "Poroto verde"
$ws.Range("H579").Value = "Magnum"
$ws.Range("I579").Value = "Primera"
$ws.Range("J579").Value = 73
$ws.Range("K579").Value = 27000
$ws.Range("L579").Value = 28000
$ws.Range("M579").Value = 27521
$ws.Range("N579").Value = "`$/malla 25 kilos"
$ws.Range("O579").Value = "Provincia de Limarí"
$ws.Range("P579").Value = 1101
$ws.Range("Q579").Value = 25
$ws.Range("R579").Value = "Hortaliza"

# --- New row 580 ---
$ws.Range("A580").Value = 3
$ws.Range("B580").Value = "Femacal de La Calera"
$ws.Range("C580").Value = "Coquimbo"
$ws.Range("D580").Value = 45075
$ws.Range("E580").Value = 5
$ws.Range("F580").Value = 100112031
$ws.Range("G580").Value = "Poroto verde"
$ws.Range("H580").Value = "Magnum"
$ws.Range("I580").Value = "Segunda"
$ws.Range("J580").Value = 38
$ws.Range("K580").Value = 21000
$ws.Range("L580").Value = 21000
$ws.Range("M580").Value = 21000
$ws.Range("N580").Value = "`$/malla 25 kilos"
$ws.Range("O580").Value = "Provincia de Limarí"
$ws.Range("P580").Value = 840
$ws.Range("Q580").Value = 25
$ws.Range("R580").Value = "Hortaliza"
